$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7281404137611389
$ws.Range("B1").Value = 3.062138319015503
$ws.Range("C1").Value = 3.798005104064941
$ws.Range("D1").Value = 1.201214194297791
$ws.Range("E1").Value = 0.944490909576416
